$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "lang_code"
$ws.Range("B1").Value = "id"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "descr"
$ws.Range("E1").Value = "is_active"

# Data rows: lang_code, id, name, descr, is_active
$data = @(
    @("eng", 10001, "Pre-Registration", "Web portal for pre-registrations", $true),
    @("eng", 10002, "Registration Client", "Desktop application for Registrations", $true),
    @("eng", 10003, "Registration Processor", "Application for post-registration process", $true),
    @("eng", 10004, "ID Authentication", "Application for third party service provider authentication", $true),
    @("eng", 10005, "ID Control", "Web portal for configuring applications", $true),
    @("eng", 10006, "Resident Portal", "Web portal for Post ID generation services", $true),
    @("fra", 10001, "Pré-inscription", "Portail Web pour les pré-inscriptions", $true),
    @("fra", 10002, "Client dinscription", "Application de bureau pour les inscriptions", $true),
    @("fra", 10003, "Processeur dinscription", "Demande de post-inscription", $true),
    @("fra", 10004, "Authentification ID", "Application pour lauthentification du fournisseur de services tiers", $true),
    @("fra", 10005, "Contrôle didentité", "Portail Web pour la configuration dapplications", $true),
    @("fra", 10006, "Portail Résident", "Portail Web pour les services de génération de post-ID", $true)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# Apply the same header-style look (bold, centered, thin box border) to the
# lang_code data column (A2:A13) as seen on the header row.
$rng = $ws.Range("A2:A13")
$rng.Font.Bold = $true
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4160
$rng.Borders.LineStyle = 1
